$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'317.68"
$ws.Range("E2").Value = "'3.83%"
$ws.Range("D3").Value = "'39.78"
$ws.Range("E3").Value = "'2.39%"
$ws.Range("D4").Value = "'5.145"
$ws.Range("E4").Value = "'0.74%"
$ws.Range("D5").Value = "'0.08213"
$ws.Range("E5").Value = "'1.58%"
$ws.Range("D6").Value = "'2.070"
$ws.Range("E6").Value = "'7.09%"
$ws.Range("D7").Value = "'8.360"
$ws.Range("E7").Value = "'4.57%"
$ws.Range("D8").Value = "'0.9406"
$ws.Range("E8").Value = "'1.01%"
$ws.Range("D9").Value = "'0.1363"
$ws.Range("E9").Value = "'-6.52%"
$ws.Range("D10").Value = "'0.1976"
$ws.Range("E10").Value = "'3.10%"
$ws.Range("D11").Value = "'0.09102"
$ws.Range("E11").Value = "'0.15%"
$ws.Range("D12").Value = "'0.03512"
$ws.Range("E12").Value = "'0.04%"
$ws.Range("E13").Value = "'0.26%"
$ws.Range("D14").Value = "'0.001417"
$ws.Range("E14").Value = "'1.92%"
$ws.Range("D15").Value = "'0.006213"
$ws.Range("E15").Value = "'6.31%"
$ws.Range("E16").Value = "'-2.26%"
$ws.Range("D17").Value = "'4.338"
$ws.Range("E17").Value = "'3.51%"
$ws.Range("D18").Value = "'3.337"
$ws.Range("E18").Value = "'-3.31%"
$ws.Range("D19").Value = "'0.3477"
$ws.Range("E19").Value = "'0.99%"
$ws.Range("D20").Value = "'0.1311"
$ws.Range("E20").Value = "'-2.82%"
$ws.Range("D21").Value = "'4.950"
$ws.Range("E21").Value = "'5.77%"
$ws.Range("E22").Value = "'1.37%"
$ws.Range("D23").Value = "'0.04356"
$ws.Range("E23").Value = "'-0.34%"
$ws.Range("D24").Value = "'0.001229"
$ws.Range("E24").Value = "'-0.73%"
$ws.Range("D25").Value = "'0.004816"
$ws.Range("E25").Value = "'12.66%"
$ws.Range("D26").Value = "'0.0001300"
$ws.Range("E26").Value = "'-0.24%"
$ws.Range("D27").Value = "'0.0004002"
$ws.Range("E27").Value = "'-10.03%"
$ws.Range("D39").Value = "'0.02246"
$ws.Range("E39").Value = "'10.37%"
$ws.Range("D40").Value = "'0.05206"
$ws.Range("E40").Value = "'2.83%"
$ws.Range("D41").Value = "'0.007766"
$ws.Range("E41").Value = "'3.33%"
$ws.Range("D42").Value = "'0.009698"
$ws.Range("E42").Value = "'-0.48%"
$ws.Range("E43").Value = "'4.59%"
$ws.Range("D44").Value = "'0.002049"
$ws.Range("E44").Value = "'-3.62%"
$ws.Range("D45").Value = "'0.009667"
$ws.Range("E45").Value = "'-2.62%"
$ws.Range("D46").Value = "'0.00006656"
$ws.Range("E46").Value = "'7.36%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.14%"
$ws.Range("D48").Value = "'0.002944"
$ws.Range("E48").Value = "'2.35%"
$ws.Range("D49").Value = "'0.001691"
$ws.Range("E49").Value = "'-6.25%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.14%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.14%"
